$d = $word.ActiveDocument

# The "Research Progress" section intro paragraph gets folded into the
# "Meeting: ..." sentence (its own Heading2 + FirstParagraph are removed as
# separate paragraphs, their text becomes a continuation of the prior run).
$mergedText = "Meeting: July 8th, 2015 in MSB 5503 ## Research Progress Over the past year I have worked on tracking a set of thirty sequence-discrete populations from Trout Bog, as they are represented by genomes assembled from metagenomes(GFMs), using 63 metagenomes from 6 of years. I called single nucleotide polymorphisms (SNPs) and looked for genes gained or lost by the populations. One genome showed a genome-wide sweep of diversity through the time series. Other genomes show evidence of past gene-sweeps, regions of statistically significant low diversity. From these results, we have propose that diversity within different co-occuring populations may be controlled by different mechanisms(recombintaiton vs. selection). I currently have a manuscript submitted to ISMEJ on the results of this analysis."

$d.Content.Find.Execute("Meeting: July 8th, 2015 in MSB 5503", $true, $false, $false, $false, $false, $true, 1, $false, $mergedText, 2) | Out-Null

# Now remove the (now-duplicate) "Research Progress" Heading2 paragraph and
# the "Over the past year ..." FirstParagraph that followed it, since that
# text now lives in the paragraph above.
$headingPara = $d.Paragraphs.Item(3)
$bodyPara = $d.Paragraphs.Item(4)
$toRemove = $d.Range($headingPara.Range.Start, $bodyPara.Range.End)
$toRemove.Delete()
